$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column N by copying the formatting from column M (same row) and
# then writing in the new 2020 data values.

$ws.Range("M4").Copy($ws.Range("N4"))
$ws.Range("N4").Value = 2020

$ws.Range("M5").Copy($ws.Range("N5"))
$ws.Range("N5").Value = 11.4

$ws.Range("M6").Copy($ws.Range("N6"))
$ws.Range("N6").Value = 14.7

$ws.Range("M7").Copy($ws.Range("N7"))
$ws.Range("N7").Value = 9

$ws.Range("M8").Copy($ws.Range("N8"))
$ws.Range("N8").Value = 10.8

$ws.Range("M9").Copy($ws.Range("N9"))
$ws.Range("N9").Value = 4.7

$ws.Range("M10").Copy($ws.Range("N10"))
$ws.Range("N10").Value = 5.0999999999999996

$ws.Range("M11").Copy($ws.Range("N11"))
$ws.Range("N11").Value = 3.4

$ws.Range("M12").Copy($ws.Range("N12"))
$ws.Range("N12").Value = 19.7

$ws.Range("M13").Copy($ws.Range("N13"))
$ws.Range("N13").Value = 18.8

$ws.Range("M14").Copy($ws.Range("N14"))
$ws.Range("N14").Value = 6.8

$ws.Range("M15").Copy($ws.Range("N15"))

$ws.Range("M16").Copy($ws.Range("N16"))
$ws.Range("N16").Value = 12.5

$ws.Range("M17").Copy($ws.Range("N17"))
$ws.Range("N17").Value = 10.7

# Update the sheet view: scroll position and active selection.
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("E2").Select()
$ws.Range("S18").Select()
